# Edit: add "PO Forecast" sheet, rename headers on existing sheets
$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3) Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- 4) Data rows: ds, PO_Forecast, yhat_lower, yhat_upper ---
$data = @(
    @(44983.99999999999,4,-65.28181367163427,74.0847763273791),
    @(44990.99999999999,5,-64.62105552672006,70.46036673757422),
    @(44997.99999999999,6,-59.56060092631959,71.53540092137062),
    @(45004.99999999999,7,-61.63569834383136,73.86807469086288),
    @(45011.99999999999,8,-59.60773313406054,73.48297612634865),
    @(45018.99999999999,8,-63.09722005794026,72.22597556135361),
    @(45025.99999999999,9,-57.27728654021232,75.66855896049529),
    @(45039.99999999999,11,-56.18797111989303,80.2253572838858),
    @(45046.99999999999,12,-55.39682218101226,78.43385609209632),
    @(45053.99999999999,13,-54.02798924953189,73.75900567936324),
    @(45088.99999999999,17,-52.17737579312664,82.81019034625507),
    @(45109.99999999999,19,-45.49904764809901,84.13563921222048),
    @(45116.99999999999,20,-48.42198239338045,86.15898256629883),
    @(45130.99999999999,22,-47.09371819143753,87.25620737860071),
    @(45137.99999999999,22,-43.02047714252994,86.35429072184765),
    @(45144.99999999999,23,-41.22751276521442,90.47717844141401),
    @(45151.99999999999,24,-41.74389874319579,95.25360401624272),
    @(45158.99999999999,25,-42.51503244873866,91.15622805058612),
    @(45165.99999999999,26,-41.02365396066121,88.00285007675129),
    @(45172.99999999999,26,-39.42874052424613,94.51086558380278),
    @(45179.99999999999,27,-40.48261188923593,93.01467670242556),
    @(45186.99999999999,28,-33.9570665493616,92.64828313783755),
    @(45193.99999999999,29,-33.63269139541328,98.39970051705647),
    @(45214.99999999999,31,-37.29502044554864,98.00832298852691),
    @(45221.99999999999,32,-32.21429094870361,107.6921133523559),
    @(45228.99999999999,33,-38.12970523791322,94.54010196235446),
    @(45235.99999999999,34,-30.81913415069404,98.26848875126881),
    @(45242.99999999999,35,-31.01217754231835,108.6347826308741),
    @(45249.99999999999,36,-30.49235340279209,102.9487808956994),
    @(45270.99999999999,38,-30.23083586470615,107.3409166031561),
    @(45277.99999999999,39,-30.35606462836303,109.3147430975068),
    @(45298.99999999999,41,-31.02179409760469,104.1576471984174),
    @(45340.99999999999,46,-24.51241379197113,116.047044328333),
    @(45382.99999999999,51,-9.184015737809762,117.0135003587168),
    @(45389.99999999999,52,-17.61786820833674,116.8314810150924),
    @(45396.99999999999,53,-17.57403782060511,118.0028512082258),
    @(45403.99999999999,54,-14.9665779676802,125.0548046005801),
    @(45410.99999999999,54,-11.95238873996949,121.5751688052979),
    @(45417.99999999999,55,-9.672182541703661,116.1979755336844),
    @(45431.99999999999,57,-9.891720196945494,123.7649168243463),
    @(45438.99999999999,58,-11.29081800858543,124.1621458594554),
    @(45445.99999999999,58,-6.636912765887504,129.4023938469066),
    @(45452.99999999999,59,-9.439973593954162,126.3139642217555),
    @(45459.99999999999,60,-4.694744908005594,124.8272430951001),
    @(45466.99999999999,61,-6.522013724947867,130.5346266668779),
    @(45473.99999999999,62,-8.756987495648501,124.7777774517368),
    @(45480.99999999999,63,-4.150890099717522,122.1038754121277),
    @(45487.99999999999,63,-3.640570599332932,127.3912473686233),
    @(45494.99999999999,64,-5.110381314748742,135.7411233796234),
    @(45543.99999999999,70,6.11513133627529,131.2550253383951),
    @(45620.99999999999,79,10.7251693322595,145.9299625420225),
    @(45627.99999999999,80,18.01614202416655,151.087284927588),
    @(45634.99999999999,81,12.09844661493914,142.700126388673),
    @(45641.99999999999,81,13.73843893004617,151.0994423346184),
    @(45648.99999999999,82,15.81345654696508,146.7344620146449),
    @(45655.99999999999,83,15.49251450104773,150.0732993435449),
    @(45662.99999999999,84,18.80795236518971,153.32156879878),
    @(45669.99999999999,85,15.65951206529355,146.9287290879855),
    @(45676.99999999999,85,18.80461498607993,151.4064603193542),
    @(45683.99999999999,86,16.81090055113378,155.5834934161426),
    @(45690.99999999999,87,21.29948982718029,150.5294107107808),
    @(45697.99999999999,88,19.26758406067586,153.9812877518465)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Apply the date/time number format to column A data cells (matches "Order Week"/"Order Month" columns)
$wsForecast.Range("A2:A63").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "PO Forecast sheet added with $($data.Count) data rows"
